# Updates the "Price" (D) and "Volume(1h)" (E) columns on the cryptos sheet
# with freshly scraped values. All of these columns hold plain text in the
# workbook (prices use "." as a thousands separator, e.g. "34.854.56", and
# volumes are padded percent strings like "  -0.83%  "), so every value is
# written/kept as text rather than letting Excel infer a numeric type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are NOT parsed as plain numbers by Excel (percent
# strings, and prices that contain more than one "." like "34.854.56").
# A direct .Value assignment already keeps these as text.
$textUpdates = @{
    "D2" = "34.854.56"
    "E2" = "  -0.83%  "
    "D3" = "1.839.79"
    "E3" = "  +1.24%  "
    "E4" = "  +0.00%  "
    "E5" = "  -0.70%  "
    "E6" = "  +0.91%  "
    "E7" = "  -0.04%  "
    "E8" = "  -4.81%  "
    "E9" = "  +0.62%  "
    "E10" = "  -0.28%  "
    "E11" = "  -1.39%  "
    "D12" = "2.104.16"
    "E12" = "  +1.09%  "
    "E13" = "  +2.14%  "
    "D14" = "1.836.47"
    "E14" = "  +0.43%  "
    "E16" = "  -0.28%  "
    "D17" = "34.855.67"
    "E17" = "  -0.74%  "
    "E18" = "  +0.06%  "
    "E19" = "  -0.86%  "
    "E20" = "  +0.35%  "
    "E21" = "  +1.97%  "
    "E22" = "  +0.57%  "
    "E23" = "  +0.23%  "
    "E24" = "  +0.01%  "
    "E25" = "  -0.56%  "
    "E26" = "  -0.71%  "
    "E27" = "  +2.46%  "
    "E28" = "  -0.49%  "
    "E29" = "  -5.77%  "
    "E30" = "  -0.04%  "
    "E31" = "  -0.41%  "
    "E32" = "  -3.21%  "
    "E33" = "  -1.42%  "
    "E34" = "  +4.51%  "
    "E35" = "  +6.67%  "
    "E36" = "  +11.45%  "
    "E37" = "  +1.56%  "
    "E39" = "  +5.25%  "
    "D40" = "1.340.55"
    "E40" = "  +2.34%  "
    "E41" = "  -0.50%  "
    "E42" = "  -0.91%  "
    "E43" = "  -1.62%  "
    "E44" = "  -2.46%  "
    "E45" = "  -0.17%  "
    "E46" = "  -0.97%  "
    "E47" = "  +1.74%  "
    "D48" = "2.018.38"
    "E48" = "  +1.06%  "
    "E49" = "  +5.26%  "
    "E50" = "  -0.03%  "
    "E51" = "  +14.45%  "
}

foreach ($cellRef in $textUpdates.Keys) {
    $ws.Range($cellRef).Value = $textUpdates[$cellRef]
}

# Cells whose new values look like plain decimal numbers (e.g. "231.01",
# "17.40", "0.0680"). Excel would otherwise silently convert these to
# numeric cells (dropping the original text formatting, e.g. trailing
# zeros), so force Text format first, then restore the cell style to
# Normal/General afterwards so no stray formatting is left on the cell.
$numericLookingUpdates = @{
    "D5" = "231.01"
    "D6" = "0.620"
    "D8" = "39.61"
    "D9" = "0.328"
    "D11" = "0.0986"
    "D15" = "0.672"
    "D16" = "4.64"
    "D18" = "69.71"
    "D20" = "240.51"
    "D21" = "12.17"
    "D24" = "2.26"
    "D25" = "171.45"
    "D26" = "7.79"
    "D28" = "17.40"
    "D34" = "1.86"
    "D37" = "0.694"
    "D38" = "91.49"
    "D39" = "1.05"
    "D42" = "14.58"
    "D43" = "2.28"
    "D46" = "6.27"
    "D49" = "0.0680"
}

foreach ($cellRef in $numericLookingUpdates.Keys) {
    $ws.Range($cellRef).NumberFormat = "@"
}
foreach ($cellRef in $numericLookingUpdates.Keys) {
    $ws.Range($cellRef).Value = $numericLookingUpdates[$cellRef]
}
foreach ($cellRef in $numericLookingUpdates.Keys) {
    $ws.Range($cellRef).NumberFormat = "General"
    $ws.Range($cellRef).Style = "Normal"
}
